# Add a "Unrelated data" sheet (Sage 50 dashboard widget export) with a
# dropdown-driving table, and rename the original export sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sage 50 journals export"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Unrelated data"

# --- Header row (shared-string insertion order: Divisional output,
# Widgets bought, Widgets sold -- matches the source workbook's table) --
$ws2.Range("A1").Value = "Date"
$ws2.Range("C1").Value = "Divisional output"
$ws2.Range("B1").Value = "Widgets bought"
$ws2.Range("D1").Value = "Widgets sold"
$ws2.Range("A1:D1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 31.5

# --- Data rows (Date, Widgets bought, Divisional output, Widgets sold) -
$data = @(
  @(41839, 5, 2000, 1),
  @(41842, 1, 3316, 0),
  @(41843, 3, 5617, 0),
  @(41849, 1, 1500, 0),
  @(41850, 0, 1500, 1),
  @(41853, 0, 2000, 0),
  @(41858, 2, 2000, 0),
  @(41860, 2, 2000, 0),
  @(41867, 2, 25468, 2),
  @(41875, 1, 43217, 1),
  @(41878, 3, 11245, 0),
  @(41881, 0, 17930, 1),
  @(41895, 2, 40799, 2),
  @(41902, 0, 20026, 0),
  @(41905, 1, 17353, 2),
  @(41909, 0, 41325, 0),
  @(41916, 3, 42713, 1),
  @(41930, 0, 29615, 8),
  @(41937, 0, 44449, 2),
  @(41946, 3, 23212, 1),
  @(41952, 1, 43476, 1),
  @(41965, 0, 31825, 0),
  @(41972, 0, 45232, 0),
  @(41976, 1, 41152, 4),
  @(41979, 0, 44716, 0),
  @(41986, 1, 41694, 1),
  @(41994, 1, 52315, 0),
  @(41999, 1, 44817, 3),
  @(42001, 0, 35436, 0),
  @(42005, 2, 45367, 3),
  @(42008, 1, 30302, 0),
  @(42014, 0, 45369, 1),
  @(42021, 1, 35973, 2),
  @(42028, 0, 22961, 0)
)

$r = 2
foreach ($row in $data) {
  $ws2.Cells.Item($r, 1).Value = $row[0]
  $ws2.Cells.Item($r, 2).Value = $row[1]
  $ws2.Cells.Item($r, 3).Value = $row[2]
  $ws2.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$ws2.Range("A2:A35").NumberFormat = "m/d/yy"
$ws2.Range("B2:B35").NumberFormat = "0"
$ws2.Range("D2:D35").NumberFormat = "0"
$ws2.Range("C2:C35").NumberFormat = "#,##0"
$ws2.Range("A2:D35").VerticalAlignment = -4108
$ws2.Range("A2:D35").WrapText = $true

# --- Column widths / styles -------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 10.7109375
$ws2.Columns.Item(2).ColumnWidth = 9.140625
$ws2.Columns.Item(3).ColumnWidth = 9.7109375
$ws2.Columns.Item(4).ColumnWidth = 8.28515625

$ws2.Range("H14").Select()
